$p = $ppt.ActivePresentation

# --- Slide 11 ("What we learned") text tweaks -----------------------------
# Shape 2 ("TextShape 2") holds four bullet paragraphs; fix a double-space
# typo in the second bullet and replace the last ("Anything else?....")
# bullet with the final wrap-up sentence, while keeping each paragraph as a
# single run (so formatting/run properties are preserved as-is).
$s11 = $p.Slides.Item(11)
$bullets = $s11.Shapes.Item(2).TextFrame.TextRange

$bullets.Paragraphs(2).Runs(1).Text = "Dependency injection is best"
$bullets.Paragraphs(4).Runs(1).Text = "Breakdown in communication results in unproductive use of time and thus communication is crucial"

# --- Slide-number field housekeeping ---------------------------------------
# Touch the slide-number placeholder on every slide's master (akin to
# reopening Insert > Header & Footer and re-applying the slide-number field)
# so the <a:fld> slide-number fields get refreshed/re-cached. This deck uses
# two different masters (slide 1 vs. the rest), and re-running this against
# an already-processed master is a harmless no-op, so we don't bother
# de-duplicating - every master's placeholder gets visited.
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $m = $p.Slides.Item($i).Master
    for ($j = 1; $j -le $m.Shapes.Count; $j++) {
        $ph = $m.Shapes.Item($j)
        if ($ph.TextFrame.HasText -and $ph.TextFrame.TextRange.Text -eq "<number>") {
            $ph.TextFrame.TextRange.InsertSlideNumber()
        }
    }
}
